# Form the consolidated report: update the "Absent" column (H) values
# for the attendance summary sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows where Absent flips from 0 -> 1
$ws.Range("H3").Value = 1
$ws.Range("H9").Value = 1
$ws.Range("H14").Value = 1
$ws.Range("H19").Value = 1

# Rows where Absent cell was blank -> now filled in with 0
$ws.Range("H6").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("H16").Value = 0
$ws.Range("H20").Value = 0
